$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("A01", "seg23.1.j.d", "GACTTTTTACGAGGCAGATAAAATTTCTC"),
    @("A02", "seg23.2.j.d", "AGGAGTAAATAAAAAACAAAGAGACACAACT"),
    @("A03", "seg23.3.j.d", "CGCTGGATTTTCGACTAAAGTTCTG"),
    @("A04", "seg23.4.j.d", "TCATCATTCTCAATCAGACAACTCGA"),
    @("A05", "seg23.5.j.d", "GATGACTTACAAATGCTTAGCCAATGA"),
    @("A06", "seg23.6.j.d", "TGTTCCTGCTGACCTTCGTTC"),
    @("A07", "seg23.7.j.d", "GATTCAGCTGTAAACTCGCCAT"),
    @("A08", "seg23.8.d", "ATACAATTTGCGTGCCAGATTTTTATCTTTC"),
    @("C01", "seg23.1.j.w", "GACTTTTTACGCGGCAGATAAAATTTCTG"),
    @("C02", "seg23.2.j.w", "AGGAGTAAATAATGAACAAAGAGACACAACC"),
    @("C03", "seg23.3.j.w", "CGCTGGATTTTTGACTAAAGTTTCT"),
    @("C04", "seg23.4.j.w", "TCATCATTCTCAATCAGACAACTCTC"),
    @("C05", "seg23.5.j.w", "CATGACTTACAAATGCTTAGCCAATGG"),
    @("C06", "seg23.6.j.w", "TGTTCCTGCTGACCTTCGTTG"),
    @("C07", "seg23.7.j.w", "GAATCAGCTGTAACTGCGCCAG"),
    @("C08", "seg23.8.w", "ATACAATTTGCGTGCCAATTTTTTATCTTTT"),
    @("E01", "seg23.1.j.c", "GCTAATCTGGAACGTGACTGGTA"),
    @("E02", "seg23.2.j.c", "TCCGCGAAGAGTTGAACGAA"),
    @("E03", "seg23.3.j.c", "GCTAGGCCCTTAACGGTTATT"),
    @("E04", "seg23.4.j.c", "CCAGAAATGAATCATGGGGTATTTGAATAT"),
    @("E05", "seg23.5.j.c", "AGCGGGAAATGACAGGCA"),
    @("E06", "seg23.6.j.c", "GCGTAATCCCGGCAATCATC"),
    @("E07", "seg23.7.j.c", "CCAGCGCGGATTGTTCCC"),
    @("E08", "seg23.8.c", "GCTGTTCGCCTCACCAAGT"),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
